# Regenerate save_data to use K instead of Strike# for wittgren_nick.xlsx
# This updates the "K" column (column G) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G), as produced by the
# regenerated save_data pipeline (K replaces the old Strike# column).
$kValues = @{
    2  = 0
    4  = 0
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 2
    18 = 1
    19 = 0
    20 = 3
    21 = 1
    22 = 2
    23 = 0
    24 = 3
    25 = 2
    26 = 3
    28 = 0
    29 = 3
    30 = 2
    31 = 1
    32 = 1
    33 = 3
    34 = 0
    35 = 1
    36 = 0
    37 = 2
    38 = 1
    39 = 0
    40 = 1
    41 = 2
    42 = 1
    43 = 0
    44 = 2
    45 = 1
    46 = 0
    47 = 2
    48 = 2
    49 = 2
    50 = 1
    51 = 0
    52 = 3
    53 = 0
    54 = 1
    55 = 1
    56 = 2
    57 = 1
    58 = 1
    59 = 1
    60 = 0
    61 = 0
    62 = 1
    63 = 1
    64 = 1
    65 = 1
    66 = 3
    67 = 1
    70 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
